# Add the initial "Hash" sort (JAVA) test results that were just collected
# on the laptop, for N = 1000, 10000 and 100000 (rows 97-99). The three
# downstream average/ratio formulas in columns I, J and K are already in
# place as shared formulas and will recalc once F:H are populated, clearing
# up the #DIV/0! errors that used to show there.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F97").Value = 12.726511
$ws.Range("G97").Value = 15.246131
$ws.Range("H97").Value = 13.645118999999999

$ws.Range("F98").Value = 86.402045999999999
$ws.Range("G98").Value = 90.485462999999996
$ws.Range("H98").Value = 83.001548999999997

$ws.Range("F99").Value = 4542.8524809999999
$ws.Range("G99").Value = 4719.5939600000002
$ws.Range("H99").Value = 4667.7834579999999

# Leave the selection where the author left off editing.
[void]$ws.Range("F100").Select()
